$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the original sheet and add the new "price_types" lookup sheet right after it
$ws.Name = "vending_request_upload_template"
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "price_types"

# Populate cell values in the exact order that mirrors the shared-string table
# of the target workbook (min, max, description_1, description_2, sales_rep,
# supply_net_number, mfg, mfg_number, issue_qty, margin, price, price_type,
# profit, customer).
$ws.Range("J1").Value = "min"
$ws.Range("K1").Value = "max"
$ws.Range("B1").Value = "description_1"
$ws.Range("C1").Value = "description_2"
$ws.Range("A1").Value = "sales_rep"
$ws.Range("D1").Value = "supply_net_number"
$ws.Range("E1").Value = "mfg"
$ws.Range("F1").Value = "mfg_number"
$ws.Range("G1").Value = "issue_qty"
$ws2.Range("A2").Value = "margin"
$ws.Range("H1").Value = "price"
$ws.Range("I1").Value = "price_type"
$ws2.Range("A1").Value = "profit"
$ws.Range("L1").Value = "customer"

# Approximate the best-fit column widths of the final template (closest
# achievable via the ColumnWidth -> OOXML width conversion used here).
$ws.Columns.Item(1).ColumnWidth = 7.830729166666666
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 16.385416666666664
$ws.Columns.Item(5).ColumnWidth = 3.276041666666666
$ws.Columns.Item(6).ColumnWidth = 10.498697916666666
$ws.Columns.Item(7).ColumnWidth = 7.721354166666666
$ws.Columns.Item(8).ColumnWidth = 4.166666666666666
$ws.Columns.Item(9).ColumnWidth = 8.721354166666666
$ws.Columns.Item(10).ColumnWidth = 3.166666666666666
$ws.Columns.Item(11).ColumnWidth = 3.608072916666666
$ws.Columns.Item(12).ColumnWidth = 7.830729166666666

# Data validation on the price_type input cell, backed by the price_types sheet
$ws.Range("I2").Validation.Add(3, 1, 1, "price_types!`$A`$1:`$A`$2") | Out-Null

# Restore the user selections / active sheet seen in the final workbook
$ws2.Range("A3").Select() | Out-Null
$ws.Range("P12").Select() | Out-Null
$ws.Activate() | Out-Null
